$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "Data" to "Summary"
$ws.Name = "Summary"

# 2. Clear the old "Micro / SMEs / MSMEs" header row (row 5) and the old
#    "Enterprises (% of total)" cell (row 6) - they get re-created further
#    down the sheet below the new "Source Type" line.
$ws.Range("B5:D5").Clear()
$ws.Range("A6").Clear()

# --- Re-assert formatting on the pre-existing, untouched cells -------------
# (the underlying engine's own style table gets rebuilt on load, so the
# original named styles must be re-applied explicitly to keep the sheet
# looking right)
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# 3. New bold+underlined sub-heading under the title
$ws.Range("A8").Value = "Source Type: Ministry of Finance/Central Bank"
$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").Font.Underline = $true

# 4. Re-created bold column headers, now on row 10
$ws.Range("B10").Value = "Micro"
$ws.Range("B10").Font.Bold = $true

$ws.Range("C10").Value = "SMEs"
$ws.Range("C10").Font.Bold = $true

$ws.Range("D10").Value = "MSMEs"
$ws.Range("D10").Font.Bold = $true

# 5. Re-created bold row label, now on row 11, plus the new data value
$ws.Range("A11").Value = "Enterprises (% of total)"
$ws.Range("A11").Font.Bold = $true

# Store "93.1" as literal text (matches the shared-string cell in the
# target file) rather than letting it be auto-converted to a number.
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "93.1"

# 6. New italic source line directly below the data
$ws.Range("A12").Value = "Source: UNICONS - Central Bank of Sudan, 2005"
$ws.Range("A12").Font.Italic = $true

# 7. New bold reference heading further down the sheet
$ws.Range("A19").Value = "UNICONS - Central Bank of Sudan"
$ws.Range("A19").Font.Bold = $true

# 8. New italic full citation line
$ws.Range("A20").Value = "UNICONS - Central Bank of Sudan, ""SITUATIONAL ANALYSIS OF THE MICROFINANCE SECTOR IN SUDAN"", 2006, p. 8-9. Available at http://www.mfu.gov.sd/sites/default/files/microsoft_word_-_situational_analysis_of_mf_sector_in_sudan-unicons.pdf"
$ws.Range("A20").Font.Italic = $true
